$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 ("Crossing") updates
$ws.Range("B3").Value = 512
$ws.Range("C3").Value = 376
$ws.Range("D3").Value = 7095
$ws.Range("E3").Value = 5638
$ws.Range("F3").Value = 243693

# Row 5 ("UnstableInterface") updates
$ws.Range("B5").Value = 662
$ws.Range("E5").Value = 6946
$ws.Range("F5").Value = 299584
